$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Region 1: "Pierw" + [_GoBack bookmark] + "sza transakcja przy użyciu "
#   -> single run "Pierwsza transakcja przy użyciu " (bookmark removed
#      here; it is re-created further down in region 3).
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Pierwsza transakcja przy użyciu ", $true, $false, $false, $false,
    $false, $true, 1, $false, "Pierwsza transakcja przy użyciu ", 2) | Out-Null

# ---------------------------------------------------------------
# Region 2: "Wykonujemy " + [spellStart]"pierwszę"[spellEnd]
#   -> "Wykonujemy pierwsz" run + "ą" run (no proofErr marks).
# ---------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute(
    "Wykonujemy pierwszę ", $true, $false, $false, $false, $false, $true,
    1, $false, "Wykonujemy pierwsz ", 2) | Out-Null

$r2b = $d.Content
$r2b.Find.Execute("Wykonujemy pierwsz ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos2 = $r2b.Start + 18   # right after "...pierwsz", before the trailing space
$pos2End = $pos2 + 1

$ins2 = $d.Range($pos2, $pos2)
$ins2.InsertAfter("ą")

# Force the newly-inserted "ą" into its own run by toggling Bold off/on.
$split2 = $d.Range($pos2, $pos2End)
$split2.Bold = 1
$split2.Bold = 0

# ---------------------------------------------------------------
# Region 3: single run "Korzystając ze np. strony<nbsp>"
#   -> four runs: "Korzystając np. " / "ze " / "str" / "ony<nbsp>"
#      with the _GoBack bookmark re-inserted between "str" and "ony".
# ---------------------------------------------------------------
$d.Content.Find.Execute(
    "Korzystając ze np. strony", $true, $false, $false, $false, $false,
    $true, 1, $false, "Korzystając np. ze strony", 2) | Out-Null

$r3 = $d.Content
$r3.Find.Execute("Korzystając np. ze strony", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$base3 = $r3.Start
$tail3 = $r3.End + 1   # include the trailing non-breaking space

$cut1 = $base3 + 22
$cut2 = $base3 + 19
$cut3 = $base3 + 16

$s1 = $d.Range($cut1, $tail3)
$s1.Bold = 1
$s1.Bold = 0

$s2 = $d.Range($cut2, $tail3)
$s2.Bold = 1
$s2.Bold = 0

$s3 = $d.Range($cut3, $tail3)
$s3.Bold = 1
$s3.Bold = 0

$bmPos3 = $base3 + 22
$bmRange3 = $d.Range($bmPos3, $bmPos3)
$d.Bookmarks.Add("_GoBack", $bmRange3) | Out-Null

Write-Output "edit complete"
